# Insert a new weekly record at row 107 (pushes existing rows 107..199 down
# to 108..200) and populate it with the new price observation, matching the
# author's commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 107, shifting rows 107-199 to 108-200.
$ws.Rows("107:107").Insert()

# Fill in the data for the newly inserted row 107.
$ws.Range("A107").Value = 11
$ws.Range("B107").Value = "Vega Monumental Concepción"
$ws.Range("C107").Value = "Bíobío"
$ws.Range("D107").Value = 44790
$ws.Range("E107").Value = 8
$ws.Range("F107").Value = 100112003
$ws.Range("G107").Value = "Ajo"
$ws.Range("H107").Value = "Chino"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 330
$ws.Range("K107").Value = 22000
$ws.Range("L107").Value = 23000
$ws.Range("M107").Value = 22545
$ws.Range("N107").Value = "$/caja 10 kilos"
$ws.Range("O107").Value = "China"
$ws.Range("P107").Value = 2254
$ws.Range("Q107").Value = 10
$ws.Range("R107").Value = "Hortaliza"
